# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K (strikeout) values for rows 2-21, replacing the previous Strike# derived values.
$kValues = @{
    2  = 3
    3  = 6
    4  = 10
    5  = 6
    6  = 3
    7  = 3
    8  = 6
    9  = 5
    10 = 4
    11 = 0
    12 = 3
    13 = 4
    14 = 5
    15 = 5
    16 = 1
    17 = 3
    18 = 4
    19 = 1
    20 = 3
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
